$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 835
$ws.Range("J17").Value = 861
$ws.Range("L17").Value = 2583
$ws.Range("N17").Value = -2919
$ws.Range("H107").Value = 320.83334
$ws.Range("I107").Value = 274
$ws.Range("J107").Value = 555
$ws.Range("K107").Value = 274
$ws.Range("L107").Value = 555
$ws.Range("M107").Value = 1646
$ws.Range("N107").Value = -4395
$ws.Range("H112").Value = 4124.25
$ws.Range("I112").Value = 1500
$ws.Range("J112").Value = 4999
$ws.Range("K112").Value = 4500
$ws.Range("L112").Value = 14997
$ws.Range("M112").Value = -3392
$ws.Range("N112").Value = -17213
$ws.Range("H113").Value = 5000
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H127").Value = 2192.4443
$ws.Range("I127").Value = 2333.4285
$ws.Range("K127").Value = 7000.2855
$ws.Range("M127").Value = -2040.2855
$ws.Range("H135").Value = 1732.5834
$ws.Range("I135").Value = 614.0909
$ws.Range("K135").Value = 5526.8181
$ws.Range("M135").Value = -2991.8181
$ws.Range("H138").Value = 5220.838
$ws.Range("J138").Value = 5677.4194
$ws.Range("L138").Value = 17032.2582
$ws.Range("N138").Value = -27312.2582

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H92").Value = 23333.334
$ws.Range("J92").Value = 23333.334
$ws.Range("L92").Value = 23333.334
$ws.Range("N92").Value = -28325.334
$ws.Range("H132").Value = 3947.4285
$ws.Range("I132").Value = 3105.3333
$ws.Range("K132").Value = 9315.999899999999
$ws.Range("M132").Value = -6785.999899999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 397.5
$ws.Range("J12").Value = 300
$ws.Range("L12").Value = 300
$ws.Range("N12").Value = -636
$ws.Range("H86").Value = 1582
$ws.Range("J86").Value = 2749
$ws.Range("L86").Value = 2749
$ws.Range("N86").Value = -4995
$ws.Range("H89").Value = 1582
$ws.Range("J89").Value = 2749
$ws.Range("L89").Value = 13745
$ws.Range("N89").Value = -24977
$ws.Range("H94").Value = 656.3333
$ws.Range("I94").Value = 664.8333
$ws.Range("J94").Value = 639.3333
$ws.Range("K94").Value = 664.8333
$ws.Range("L94").Value = 639.3333
$ws.Range("M94").Value = -213.8333
$ws.Range("N94").Value = -1541.3333
$ws.Range("H99").Value = 2296.5557
$ws.Range("I99").Value = 2333.75
$ws.Range("J99").Value = 1999
$ws.Range("K99").Value = 2333.75
$ws.Range("L99").Value = 1999
$ws.Range("M99").Value = -835.75
$ws.Range("N99").Value = -4995
$ws.Range("H107").Value = 468.65216
$ws.Range("I107").Value = 212.8
$ws.Range("J107").Value = 948.375
$ws.Range("K107").Value = 212.8
$ws.Range("L107").Value = 948.375
$ws.Range("M107").Value = 1707.2
$ws.Range("N107").Value = -4788.375
$ws.Range("H134").Value = 2183.1765
$ws.Range("I134").Value = 2218.75
$ws.Range("J134").Value = 1614
$ws.Range("K134").Value = 6656.25
$ws.Range("L134").Value = 4842
$ws.Range("M134").Value = -4121.25
$ws.Range("N134").Value = -9912

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2332.5
$ws.Range("I31").Value = 2332.5
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 2332.5
$ws.Range("L31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -2037.5
$ws.Range("H34").Value = 2332.5
$ws.Range("I34").Value = 2332.5
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 2332.5
$ws.Range("L34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -2130.5
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()
$ws.Range("H94").Value = 3344.6667
$ws.Range("I94").Value = 2921.625
$ws.Range("J94").Value = 3828.1428
$ws.Range("K94").Value = 2921.625
$ws.Range("L94").Value = 3828.1428
$ws.Range("M94").Value = -2470.625
$ws.Range("N94").Value = -4730.1428
$ws.Range("H99").Value = 2249.5
$ws.Range("I99").Value = 2100
$ws.Range("J99").Value = 2399
$ws.Range("K99").Value = 2100
$ws.Range("L99").Value = 2399
$ws.Range("M99").Value = -602
$ws.Range("N99").Value = -5395
$ws.Range("H107").Value = 349.75
$ws.Range("I107").Value = 299.5
$ws.Range("J107").Value = 400
$ws.Range("K107").Value = 299.5
$ws.Range("L107").Value = 400
$ws.Range("M107").Value = 1620.5
$ws.Range("N107").Value = -4240
$ws.Range("H126").Value = 2249.5
$ws.Range("I126").Value = 2100
$ws.Range("J126").Value = 2399
$ws.Range("K126").Value = 6300
$ws.Range("L126").Value = 7197
$ws.Range("M126").Value = -3830
$ws.Range("N126").Value = -12137
$ws.Range("H134").Value = 2419.1428
$ws.Range("I134").Value = 2452.0833
$ws.Range("J134").Value = 2221.5
$ws.Range("K134").Value = 7356.249899999999
$ws.Range("L134").Value = 6664.5
$ws.Range("M134").Value = -4821.249899999999
$ws.Range("N134").Value = -11734.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 5498.857
$ws.Range("I121").Value = 695
$ws.Range("J121").Value = 6299.5
$ws.Range("K121").Value = 2085
$ws.Range("L121").Value = 18898.5
$ws.Range("M121").Value = -775
$ws.Range("N121").Value = -21518.5
$ws.Range("H129").Value = 2636.2
$ws.Range("J129").Value = 3397
$ws.Range("L129").Value = 10191
$ws.Range("N129").Value = -20191

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5519.8
$ws.Range("I80").Value = 2149.5
$ws.Range("K80").Value = 2149.5
$ws.Range("M80").Value = -1151.5
$ws.Range("H83").Value = 5519.8
$ws.Range("I83").Value = 2149.5
$ws.Range("K83").Value = 10747.5
$ws.Range("M83").Value = -5755.5
$ws.Range("H126").Value = 1999.5714
$ws.Range("I126").Value = 1999.5
$ws.Range("K126").Value = 5998.5
$ws.Range("M126").Value = -3528.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 556.25
$ws.Range("I16").Value = 556.25
$ws.Range("K16").Value = 556.25
$ws.Range("M16").Value = -386.25
$ws.Range("H22").Value = 865.7273
$ws.Range("J22").Value = 931.6667
$ws.Range("L22").Value = 931.6667
$ws.Range("N22").Value = -1521.6667
$ws.Range("H27").Value = 865.7273
$ws.Range("J27").Value = 931.6667
$ws.Range("L27").Value = 931.6667
$ws.Range("N27").Value = -1145.6667
$ws.Range("H40").Value = 5999.6665
$ws.Range("I40").Value = 5999
$ws.Range("K40").Value = 5999
$ws.Range("M40").Value = -5863
$ws.Range("H46").Value = 2055.3333
$ws.Range("I46").Value = 2055.3333
$ws.Range("K46").Value = 2055.3333
$ws.Range("M46").Value = -1867.3333
$ws.Range("H82").Value = 1666.3334
$ws.Range("I82").Value = 1749.5
$ws.Range("K82").Value = 1749.5
$ws.Range("M82").Value = -1388.5
$ws.Range("H85").Value = 1666.3334
$ws.Range("I85").Value = 1749.5
$ws.Range("K85").Value = 1749.5
$ws.Range("M85").Value = -501.5
$ws.Range("H100").Value = 3408.3076
$ws.Range("I100").Value = 3381
$ws.Range("J100").Value = 3499.3333
$ws.Range("K100").Value = 3381
$ws.Range("L100").Value = 3499.3333
$ws.Range("M100").Value = -2840
$ws.Range("N100").Value = -4581.3333
$ws.Range("H105").Value = 40000
$ws.Range("J105").Value = 40000
$ws.Range("L105").Value = 40000
$ws.Range("N105").Value = -46988
$ws.Range("H132").Value = 7480.636
$ws.Range("I132").Value = 7480.636
$ws.Range("K132").Value = 22441.908
$ws.Range("M132").Value = -19911.908
$ws.Range("H136").Value = 4140.6
$ws.Range("I136").Value = 4174.5
$ws.Range("J136").Value = 4005
$ws.Range("K136").Value = 12523.5
$ws.Range("L136").Value = 12015
$ws.Range("M136").Value = -9973.5
$ws.Range("N136").Value = -17115

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("L12").ClearContents()
$ws.Range("N12").Value = 0
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("H122").Value = 1974.625
$ws.Range("I122").Value = 1876
$ws.Range("K122").Value = 5628
$ws.Range("M122").Value = -3178
$ws.Range("H132").Value = 2775.1428
$ws.Range("I132").Value = 2775.1428
$ws.Range("K132").Value = 8325.428400000001
$ws.Range("M132").Value = -5795.428400000001
